# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stock) worksheet, populating header row 1 and data rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# ---- Header row (row 1): H1=date, I1=legislator_name, J1=legislator_id ----
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Match the header formatting used by the other header cells (bold + border)
$ws.Range("B1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Data rows (rows 2-3) ----
# The date must stay a literal text string ("2011-11-24") rather than be
# auto-converted into a date serial number. Enter it as a text formula
# (="2011-11-24") and then flatten it in-place to a static value via
# copy/paste-values, so the result is plain text with no left-over
# formula and no new number-format/style being introduced.
$ws.Range("H2").Formula = "=""2011-11-24"""
$ws.Range("H3").Formula = "=""2011-11-24"""
$ws.Range("H2:H3").Copy()
$ws.Range("H2:H3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("I2").Value = "廖國棟"
$ws.Range("I3").Value = "廖國棟"
$ws.Range("J2").Value = 962
$ws.Range("J3").Value = 962

# Match the data-row formatting used by the other data cells in those rows
$ws.Range("B2").Copy()
$ws.Range("H2:J2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B3").Copy()
$ws.Range("H3:J3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
